$d = $word.ActiveDocument

# 1. Update the imgur link in the "More Examples" spoiler section
$d.Content.Find.Execute("https://imgur.com/l22kHkv.png", $true, $false, $false, $false, $false,
                         $true, 1, $false, "https://imgur.com/QlsiThF.png", 2) | Out-Null

# 2. Merge the two credit paragraphs into one, removing the extra paragraph
$d.Content.Find.Execute("Alex, David, and co. for making Starsector" + [char]13 + "The Starsector community!", $true, $false, $false, $false, $false,
                         $true, 1, $false, "The Starsector devs (Alex, David, and co.) and the Starsector community!", 2) | Out-Null

# 3. Wrap "Age of Tribulation" with a darkgreen color tag
$d.Content.Find.Execute("[b]Age of Tribulation[/b]", $true, $false, $false, $false, $false,
                         $true, 1, $false, "[b][color=darkgreen]Age of Tribulation[/color][/b]", 2) | Out-Null

# 4. Add a new empty paragraph after "[/center]" (before the section break)
$p = $d.Paragraphs.Last
$p.Range.InsertParagraphAfter()
